$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 14 and row 15 coin name/link (Polkadot <-> WrappedEther)
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# Update Price (D) and Volume(1h) (E) columns; force D as text so values
# such as "1.004" are not reinterpreted as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.318.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.090.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09332"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.692"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.901"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.042.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001156"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06681"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.356"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.321.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.296"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.518"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.667"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.220"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.686"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.854"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02628"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06764"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6991"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.347"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2213"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6809"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.341"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.369"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.636"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000350"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.214"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.26%  "
